$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, matching the existing header formatting (same style as the rest of row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill in team record values for every data row (2-51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 81  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 80  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
